# chore: adapt column header formatting to respective input file names
#
# Renames the header row (row 1) suffixes from "_old"/"_new" to
# "_FV2310"/"_FV2404" respectively, wraps the data range in an Excel
# Table ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<Name>_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2310"
}

# Column K (11): "diff" (unchanged)
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U (12-21): "<Name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2404"
}

# Wrap the used range (A1:U68) in a native Excel Table ("Table1")
$lastRow = 68
$lastCol = 21
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# Freeze the header row (select the first cell below the header, then
# turn on FreezePanes so the resulting state is a true "frozen" pane).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A1").Select()
